$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update input values
$ws.Range("A2").Value = 12144
$ws.Range("B2").Value = 16832
$ws.Range("C3").Value = 98

# Recalculate workbook so dependent formulas pick up the new values
$excel.Calculate()

# Update the selected cell/range shown in the saved worksheet view
$ws.Activate()
$ws.Range("C4").Select()
